# The source workbook tracks weekly "Poroto verde" price records. A new
# weekly record is inserted at the top of the data block (row 9, right
# after the fixed header/lead-in rows 1-8), pushing all the existing
# records (old rows 9-102) down by one row (to 10-103).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 9 - shifts rows 9:102 down to 10:103
# and extends the used range to A1:R103 (matches the dimension change
# in the diff: A1:R102 -> A1:R103).
$ws.Rows("9:9").Insert()

# Populate the newly inserted row 9 with the new weekly record.
$ws.Range("A9").Value2 = 1
$ws.Range("B9").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C9").Value2 = "Arica y Parinacota"
$ws.Range("D9").Value2 = 45230
$ws.Range("E9").Value2 = 15
$ws.Range("F9").Value2 = 100112031
$ws.Range("G9").Value2 = "Poroto verde"
$ws.Range("H9").Value2 = "Sin especificar"
$ws.Range("I9").Value2 = "Primera"
$ws.Range("J9").Value2 = 750
$ws.Range("K9").Value2 = 950
$ws.Range("L9").Value2 = 1000
$ws.Range("M9").Value2 = 973
$ws.Range("N9").Value2 = "$/kilo"
$ws.Range("O9").Value2 = "Región de Arica y Parinacota"
$ws.Range("P9").Value2 = 973
$ws.Range("Q9").Value2 = 1
$ws.Range("R9").Value2 = "Hortaliza"
